$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1519516
$ws.Range("J17").Value = 1519516
$ws.Range("L17").Value = 4558548
$ws.Range("N17").Value = -4558884
$ws.Range("H19").Value = 1018.11536
$ws.Range("I19").Value = 720.0526
$ws.Range("K19").Value = 720.0526
$ws.Range("M19").Value = -545.0526
$ws.Range("H112").Value = 1214.9
$ws.Range("I112").Value = 875.8
$ws.Range("J112").Value = 1282.72
$ws.Range("K112").Value = 2627.4
$ws.Range("L112").Value = 3848.16
$ws.Range("M112").Value = -1519.4
$ws.Range("N112").Value = -6064.16
$ws.Range("H130").Value = 18470
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H137").Value = 1148.5714
$ws.Range("I137").Value = 1126.75
$ws.Range("J137").Value = 1177.6666
$ws.Range("K137").Value = 3380.25
$ws.Range("L137").Value = 3532.9998
$ws.Range("M137").Value = -830.25
$ws.Range("N137").Value = -8632.9998
$ws.Range("H138").Value = 3789.48
$ws.Range("I138").Value = 2488.7317
$ws.Range("J138").Value = 4693.3896
$ws.Range("K138").Value = 7466.195099999999
$ws.Range("L138").Value = 14080.1688
$ws.Range("M138").Value = -2326.195099999999
$ws.Range("N138").Value = -24360.1688
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11760.78
$ws.Range("I32").Value = 9553.724
$ws.Range("K32").Value = 9553.724
$ws.Range("M32").Value = -9266.724
$ws.Range("H74").Value = 1436.0735
$ws.Range("I74").Value = 908.12964
$ws.Range("J74").Value = 3472.4285
$ws.Range("K74").Value = 908.12964
$ws.Range("L74").Value = 3472.4285
$ws.Range("M74").Value = -34.12963999999999
$ws.Range("N74").Value = -5220.4285
$ws.Range("H77").Value = 1436.0735
$ws.Range("I77").Value = 908.12964
$ws.Range("J77").Value = 3472.4285
$ws.Range("K77").Value = 4540.6482
$ws.Range("L77").Value = 17362.1425
$ws.Range("M77").Value = -172.6481999999996
$ws.Range("N77").Value = -26098.1425
$ws.Range("H122").Value = 2368.111
$ws.Range("I122").Value = 1401
$ws.Range("K122").Value = 4203
$ws.Range("M122").Value = -1753
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6683.391
$ws.Range("I20").Value = 10832.583
$ws.Range("J20").Value = 2157
$ws.Range("K20").Value = 10832.583
$ws.Range("L20").Value = 2157
$ws.Range("M20").Value = -10585.583
$ws.Range("N20").Value = -2651
$ws.Range("H22").Value = 298.54544
$ws.Range("I22").Value = 217
$ws.Range("J22").Value = 516
$ws.Range("K22").Value = 217
$ws.Range("L22").Value = 516
$ws.Range("M22").Value = -44
$ws.Range("N22").Value = -862
$ws.Range("H134").Value = 1678.6
$ws.Range("I134").Value = 1307.1794
$ws.Range("J134").Value = 2583.9375
$ws.Range("K134").Value = 3921.5382
$ws.Range("L134").Value = 7751.8125
$ws.Range("M134").Value = -1386.5382
$ws.Range("N134").Value = -12821.8125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11114322
$ws.Range("I86").Value = 17244374
$ws.Range("J86").Value = 3602.5625
$ws.Range("K86").Value = 17244374
$ws.Range("L86").Value = 3602.5625
$ws.Range("M86").Value = -17243251
$ws.Range("N86").Value = -5848.5625
$ws.Range("H89").Value = 11114322
$ws.Range("I89").Value = 17244374
$ws.Range("J89").Value = 3602.5625
$ws.Range("K89").Value = 86221870
$ws.Range("L89").Value = 18012.8125
$ws.Range("M89").Value = -86216254
$ws.Range("N89").Value = -29244.8125
$ws.Range("H122").Value = 847928.9399999999
$ws.Range("I122").Value = 143908.86
$ws.Range("J122").Value = 1669285.6
$ws.Range("K122").Value = 431726.58
$ws.Range("L122").Value = 5007856.800000001
$ws.Range("M122").Value = -429276.58
$ws.Range("N122").Value = -5012756.800000001
$ws.Range("H133").Value = 24800
$ws.Range("J133").Value = 24800
$ws.Range("L133").Value = 24800
$ws.Range("N133").Value = -29860
$ws.Range("H134").Value = 3078.96
$ws.Range("I134").Value = 3103
$ws.Range("K134").Value = 9309
$ws.Range("M134").Value = -6774
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3060.8518
$ws.Range("I80").Value = 3001
$ws.Range("J80").Value = 3065.64
$ws.Range("K80").Value = 9003
$ws.Range("L80").Value = 9196.92
$ws.Range("M80").Value = -8067
$ws.Range("N80").Value = -11068.92
$ws.Range("H83").Value = 3060.8518
$ws.Range("I83").Value = 3001
$ws.Range("J83").Value = 3065.64
$ws.Range("K83").Value = 27009
$ws.Range("L83").Value = 27590.76
$ws.Range("M83").Value = -22329
$ws.Range("N83").Value = -36950.75999999999
$ws.Range("H113").Value = 1313.8422
$ws.Range("I113").Value = 1873.5
$ws.Range("J113").Value = 906.8182
$ws.Range("K113").Value = 5620.5
$ws.Range("L113").Value = 2720.4546
$ws.Range("M113").Value = -3450.5
$ws.Range("N113").Value = -7060.4546
$ws.Range("H122").Value = 971.2
$ws.Range("I122").Value = 587.72
$ws.Range("J122").Value = 1929.9
$ws.Range("K122").Value = 5289.48
$ws.Range("L122").Value = 17369.1
$ws.Range("M122").Value = -2839.48
$ws.Range("N122").Value = -22269.1
$ws.Range("H123").Value = 1411.8182
$ws.Range("I123").Value = 565
$ws.Range("K123").Value = 1695
$ws.Range("M123").Value = 755
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 67893.60000000001
$ws.Range("I122").Value = 97740.80499999999
$ws.Range("J122").Value = 3224.6667
$ws.Range("K122").Value = 293222.415
$ws.Range("L122").Value = 9674.000100000001
$ws.Range("M122").Value = -290772.415
$ws.Range("N122").Value = -14574.0001
$ws.Range("H126").Value = 2400.077
$ws.Range("I126").Value = 2317.3333
$ws.Range("K126").Value = 6951.999899999999
$ws.Range("M126").Value = -4481.999899999999
$ws.Range("H132").Value = 2745.353
$ws.Range("I132").Value = 2444.2173
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 7332.651899999999
$ws.Range("L132").Value = 10125
$ws.Range("M132").Value = -4802.651899999999
$ws.Range("N132").Value = -15185
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 55556896
$ws.Range("J122").Value = 2680
$ws.Range("L122").Value = 8040
$ws.Range("N122").Value = -12940
$ws.Range("H132").Value = 7854.9033
$ws.Range("I132").Value = 2787.8262
$ws.Range("K132").Value = 8363.4786
$ws.Range("M132").Value = -5833.4786
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 113752.22
$ws.Range("I122").Value = 169328.33
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 507984.99
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -505534.99
$ws.Range("N122").Value = -12700
$ws.Range("H132").Value = 1502.3623
$ws.Range("I132").Value = 1248.5869
$ws.Range("J132").Value = 2009.9131
$ws.Range("K132").Value = 3745.7607
$ws.Range("L132").Value = 6029.7393
$ws.Range("M132").Value = -1215.7607
$ws.Range("N132").Value = -11089.7393
$ws.Range("H136").Value = 4865.2812
$ws.Range("I136").Value = 3247.4375
$ws.Range("J136").Value = 6483.125
$ws.Range("K136").Value = 9742.3125
$ws.Range("L136").Value = 19449.375
$ws.Range("M136").Value = -7192.3125
$ws.Range("N136").Value = -24549.375
